$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: "Something" values become real dates (7/4/2023), keep date formatting ---
$lastUpdated = Get-Date -Year 2023 -Month 7 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Range("B2").Value = $lastUpdated.Date
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Value = $lastUpdated.Date
$ws.Range("B3").NumberFormat = "mm-dd-yy"

# --- Column C: "Thifafa" values become booleans ---
$ws.Range("C2").Value = $false
$ws.Range("C3").Value = $false

# --- New column E: Assignee ---
$ws.Range("E1").Value = "Assignee"
$ws.Range("E2").Value = "Rohit"
$ws.Range("E3").Value = "Rujuta"

# --- New column F: Last updated date (header styled w/ custom font) ---
$ws.Range("F1").Value = "Last updated date"
$ws.Range("F1").Font.Color = 7901646
$ws.Range("F1").Font.Name = "Menlo"
$ws.Range("F2").Value = "fsda"
$ws.Range("F3").Value = "fsda"

# --- Selection moves to the newly added F3 cell ---
$ws.Range("F3").Select() | Out-Null
